# Refresh the cryptos price/volume snapshot (rows 2-51, columns D "Price"
# and E "Volume(1h)"). Every D-column value is force-entered as text via a
# leading apostrophe (otherwise strings like "1.005" or "99.62" would be
# auto-parsed as numbers) and the cell style is reset back to "Normal"
# immediately after so no extra "Text" number-format style gets attached -
# the cells stay plain inline/shared strings, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.222.54'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.84%  '

$ws.Range("D3").Value = '''1.660.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.89%  '

$ws.Range("D4").Value = '''1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.35%  '

$ws.Range("D5").Value = '''218.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '

$ws.Range("D6").Value = '''0.5221'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.61%  '

$ws.Range("E7").Value = '  +0.34%  '

$ws.Range("D8").Value = '''0.2667'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("D9").Value = '''0.06322'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.06%  '

$ws.Range("D10").Value = '''21.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.25%  '

$ws.Range("D11").Value = '''0.07723'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.03%  '

$ws.Range("D12").Value = '''1.671.67'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.36%  '

$ws.Range("D13").Value = '''4.425'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.61%  '

$ws.Range("D14").Value = '''1.889.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.79%  '

$ws.Range("D15").Value = '''0.5458'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.81%  '

$ws.Range("D16").Value = '''0.0₅8221'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.50%  '

$ws.Range("D17").Value = '''64.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.13%  '

$ws.Range("D18").Value = '''26.258.12'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.86%  '

$ws.Range("E19").Value = '  +0.36%  '

$ws.Range("D20").Value = '''4.654'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.45%  '

$ws.Range("D21").Value = '''193.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.65%  '

$ws.Range("D22").Value = '''10.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.22%  '

$ws.Range("D23").Value = '''6.068'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.26%  '

$ws.Range("D24").Value = '''1.008'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.51%  '

$ws.Range("D25").Value = '''138.66'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.74%  '

$ws.Range("D26").Value = '''0.1240'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.84%  '

$ws.Range("D27").Value = '''7.224'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.67%  '

$ws.Range("D28").Value = '''16.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.87%  '

$ws.Range("D29").Value = '''1.402'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.66%  '

$ws.Range("D30").Value = '''0.05991'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.44%  '

$ws.Range("D31").Value = '''1.281'
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").Value = '''3.628'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.35%  '

$ws.Range("D33").Value = '''3.308'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.01%  '

$ws.Range("D34").Value = '''1.630'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.47%  '

$ws.Range("D35").Value = '''0.9787'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.64%  '

$ws.Range("E36").Value = '  -0.35%  '

$ws.Range("D37").Value = '''2.783'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.10%  '

$ws.Range("D38").Value = '''0.5881'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.82%  '

$ws.Range("D39").Value = '''0.01592'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.76%  '

$ws.Range("D40").Value = '''5.941'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.50%  '

$ws.Range("D41").Value = '''0.8609'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.14%  '

$ws.Range("E42").Value = '  +0.33%  '

$ws.Range("D43").Value = '''1.029.77'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.86%  '

$ws.Range("D44").Value = '''99.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.39%  '

$ws.Range("D45").Value = '''1.802.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.26%  '

$ws.Range("E46").Value = '  +3.67%  '

$ws.Range("D47").Value = '''57.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.18%  '

$ws.Range("E48").Value = '  -0.21%  '

$ws.Range("D49").Value = '''8.104'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.39%  '

$ws.Range("D50").Value = '''0.05183'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.44%  '

$ws.Range("D51").Value = '''1.471'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.01%  '
